$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 412-413; existing rows 412..526 shift down to 414..528
$ws.Rows("412:413").Insert()

# New row 412: Cebollín, 2023-01-06, Provincia de Cautín
$ws.Range("A412").Value = 10
$ws.Range("B412").Value = "Vega Modelo de Temuco"
$ws.Range("C412").Value = "La Araucanía"
$ws.Range("D412").Value = 44932
$ws.Range("E412").Value = 9
$ws.Range("F412").Value = 100112037
$ws.Range("G412").Value = "Cebollín"
$ws.Range("H412").Value = "Sin especificar"
$ws.Range("I412").Value = "Primera"
$ws.Range("J412").Value = 40
$ws.Range("K412").Value = 8000
$ws.Range("L412").Value = 8000
$ws.Range("M412").Value = 8000
$ws.Range("N412").Value = "`$/docena de paquetes"
$ws.Range("O412").Value = "Provincia de Cautín"
$ws.Range("P412").Value = 667
$ws.Range("Q412").Value = 12
$ws.Range("R412").Value = "Hortaliza"

# New row 413: Cebollín, 2023-01-06, Región de O'Higgins
$ws.Range("A413").Value = 10
$ws.Range("B413").Value = "Vega Modelo de Temuco"
$ws.Range("C413").Value = "La Araucanía"
$ws.Range("D413").Value = 44932
$ws.Range("E413").Value = 9
$ws.Range("F413").Value = 100112037
$ws.Range("G413").Value = "Cebollín"
$ws.Range("H413").Value = "Sin especificar"
$ws.Range("I413").Value = "Primera"
$ws.Range("J413").Value = 80
$ws.Range("K413").Value = 8000
$ws.Range("L413").Value = 8000
$ws.Range("M413").Value = 8000
$ws.Range("N413").Value = "`$/docena de paquetes"
$ws.Range("O413").Value = "Región de O'Higgins"
$ws.Range("P413").Value = 667
$ws.Range("Q413").Value = 12
$ws.Range("R413").Value = "Hortaliza"
